# Add "2022-Q3" sheet + data, matching the author's "feat: add 2022-Q3 data" commit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a range of cells to hold literal TEXT values without minting
# a new cell style (NumberFormat="@" / quotePrefix both mint a new style xf,
# which the target XML does not have). Trick: write a `="literal"` formula
# (t="str", no style change) then Copy + PasteSpecial(values-only) to freeze
# it into a plain string cell - still no style change.
# ---------------------------------------------------------------------------
function Set-TextCell {
    param($cell, [string]$text)
    $escaped = $text.Replace("""", """""")
    $cell.Formula = "=""" + $escaped + """"
}

# ===========================================================================
# 1) Insert the new "2022-Q3" worksheet right before "2022-Q2", by cloning
#    "2022-Q2" (so it starts with identical headers/styles/column layout)
#    and then overwriting its data.
# ===========================================================================
$refSheet = $wb.Worksheets.Item("2022-Q2")
$refSheet.Copy($refSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Reference sheet had 6 data rows (rows 2-7); the new sheet only needs 5
# (rows 2-6), so drop the extra row.
$newSheet.Rows.Item(7).Delete()

# New fund-holdings data for 2022-Q3.
$q3rows = @(
    @{ A = 0; B = "002121"; C = "广发沪港深新起点股票A";                     D = "25.97"; E = "85.60"; F = "4.40"; G = "1.1427"; H = 9 },
    @{ A = 1; B = "159636"; C = "工银瑞信国证港股通科技ETF";                 D = "8.42";  E = "98.06"; F = "3.36"; G = "0.2829"; H = 9 },
    @{ A = 2; B = "501021"; C = "华宝标普香港上市中国中小盘指数（LOF）A";    D = "4.19";  E = "92.99"; F = "1.61"; G = "0.0675"; H = 10 },
    @{ A = 3; B = "010024"; C = "广发沪港深新起点股票C";                     D = "0.49";  E = "85.60"; F = "4.40"; G = "0.0216"; H = 9 },
    @{ A = 4; B = "006127"; C = "华宝标普香港上市中国中小盘指数（LOF）C";    D = "0.24";  E = "92.99"; F = "1.61"; G = "0.0039"; H = 10 }
)

for ($i = 0; $i -lt $q3rows.Count; $i++) {
    $r = $i + 2
    $row = $q3rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $newSheet.Cells.Item($r, 2) $row.B
    Set-TextCell $newSheet.Cells.Item($r, 3) $row.C
    Set-TextCell $newSheet.Cells.Item($r, 4) $row.D
    Set-TextCell $newSheet.Cells.Item($r, 5) $row.E
    Set-TextCell $newSheet.Cells.Item($r, 6) $row.F
    Set-TextCell $newSheet.Cells.Item($r, 7) $row.G
    $newSheet.Cells.Item($r, 8).Value = $row.H
}

# Freeze the `="literal"` formulas above into plain static text values
# (no style change - unlike NumberFormat="@").
$textRange = $newSheet.Range("B2:G6")
$textRange.Copy()
$newSheet.Range("B2").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ===========================================================================
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 at the
#    top of the data (row 2), shifting the existing quarters down by one.
# ===========================================================================
$summary = $wb.Worksheets.Item("总计")

$existing = @()
for ($r = 2; $r -le 8; $r++) {
    $existing += ,@{
        B = $summary.Cells.Item($r, 2).Value()
        C = $summary.Cells.Item($r, 3).Value()
        D = $summary.Cells.Item($r, 4).Value()
    }
}

# New row 9 needs the same "index" column style (s="2") as the rest of the
# A column; copy it down from row 8 before that row's content gets moved.
$summary.Cells.Item(8, 1).Copy($summary.Cells.Item(9, 1))

for ($i = $existing.Count - 1; $i -ge 0; $i--) {
    $r = $i + 3
    $summary.Cells.Item($r, 1).Value = $i + 1
    $summary.Cells.Item($r, 2).Value = $existing[$i].B
    $summary.Cells.Item($r, 3).Value = $existing[$i].C
    $summary.Cells.Item($r, 4).Value = $existing[$i].D
}

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 1.52
